$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr_B_E = New-Object 'object[,]' 24,4
$arr_B_E[0,0] = 10.06602685295606
$arr_B_E[0,1] = 5.109432809453438
$arr_B_E[0,2] = 14.92169772828122
$arr_B_E[0,3] = 16.33195972981394
$arr_B_E[1,0] = 9.785404360535562
$arr_B_E[1,1] = 4.905436515945269
$arr_B_E[1,2] = 14.85612703063621
$arr_B_E[1,3] = 16.26633236098847
$arr_B_E[2,0] = 9.61047425522723
$arr_B_E[2,1] = 4.776312103355453
$arr_B_E[2,2] = 14.81905760371421
$arr_B_E[2,3] = 16.22956633925089
$arr_B_E[3,0] = 9.538637899006394
$arr_B_E[3,1] = 4.722797044463077
$arr_B_E[3,2] = 14.80476619663112
$arr_B_E[3,3] = 16.21548317490108
$arr_B_E[4,0] = 9.526679583547793
$arr_B_E[4,1] = 4.713859241059634
$arr_B_E[4,2] = 14.80244266020945
$arr_B_E[4,3] = 16.21319931280759
$arr_B_E[5,0] = 9.609507526023203
$arr_B_E[5,1] = 4.775593898906147
$arr_B_E[5,2] = 14.81886155082721
$arr_B_E[5,3] = 16.22937275286847
$arr_B_E[6,0] = 9.969877094980554
$arr_B_E[6,1] = 5.039941698424089
$arr_B_E[6,2] = 14.89843426334538
$arr_B_E[6,3] = 16.3086062833461
$arr_B_E[7,0] = 10.6512498792453
$arr_B_E[7,1] = 5.524667719817321
$arr_B_E[7,2] = 15.07921157837997
$arr_B_E[7,3] = 16.4914189090016
$arr_B_E[8,0] = 11.13095006659247
$arr_B_E[8,1] = 5.856892643377853
$arr_B_E[8,2] = 15.22625688237415
$arr_B_E[8,3] = 16.641624626044
$arr_B_E[9,0] = 11.34362987781271
$arr_B_E[9,1] = 6.002272324003365
$arr_B_E[9,2] = 15.29604068724791
$arr_B_E[9,3] = 16.71321018233084
$arr_B_E[10,0] = 11.4232956730691
$arr_B_E[10,1] = 6.056458059529687
$arr_B_E[10,2] = 15.32286375353901
$arr_B_E[10,3] = 16.74076794012682
$arr_B_E[11,0] = 11.40617809495041
$arr_B_E[11,1] = 6.044827292689463
$arr_B_E[11,2] = 15.31706955035028
$arr_B_E[11,3] = 16.73481317719374
$arr_B_E[12,0] = 11.35020190001818
$arr_B_E[12,1] = 6.006747786567089
$arr_B_E[12,2] = 15.29823957956263
$arr_B_E[12,3] = 16.71546846782865
$arr_B_E[13,0] = 11.31579927286362
$arr_B_E[13,1] = 5.983309086663457
$arr_B_E[13,2] = 15.28675687986459
$arr_B_E[13,3] = 16.70367730103238
$arr_B_E[14,0] = 11.11693291854563
$arr_B_E[14,1] = 5.847272514359705
$arr_B_E[14,2] = 15.22175298470351
$arr_B_E[14,3] = 16.63701031229346
$arr_B_E[15,0] = 10.99346069674481
$arr_B_E[15,1] = 5.762315621594834
$arr_B_E[15,2] = 15.18260306578513
$arr_B_E[15,3] = 16.59693345167268
$arr_B_E[16,0] = 10.92192547867525
$arr_B_E[16,1] = 5.712911310676579
$arr_B_E[16,2] = 15.16035844319666
$arr_B_E[16,3] = 16.57418994624345
$arr_B_E[17,0] = 10.89761844745666
$arr_B_E[17,1] = 5.696092520059337
$arr_B_E[17,2] = 15.15287429855959
$arr_B_E[17,3] = 16.56654277115906
$arr_B_E[18,0] = 11.00665862657133
$arr_B_E[18,1] = 5.77141555229167
$arr_B_E[18,2] = 15.18674247945163
$arr_B_E[18,3] = 16.60116799242546
$arr_B_E[19,0] = 11.36666767846167
$arr_B_E[19,1] = 6.01795646034008
$arr_B_E[19,2] = 15.3037597620242
$arr_B_E[19,3] = 16.72113841915817
$arr_B_E[20,0] = 11.59684519135764
$arr_B_E[20,1] = 6.174020753763321
$arr_B_E[20,2] = 15.38254496110318
$arr_B_E[20,3] = 16.80215910248856
$arr_B_E[21,0] = 11.47448557480391
$arr_B_E[21,1] = 6.091201289009415
$arr_B_E[21,2] = 15.34029093481002
$arr_B_E[21,3] = 16.7586840132736
$arr_B_E[22,0] = 11.00069354794777
$arr_B_E[22,1] = 5.767303219402541
$arr_B_E[22,2] = 15.18487022934752
$arr_B_E[22,3] = 16.59925262938201
$arr_B_E[23,0] = 10.47022424396229
$arr_B_E[23,1] = 5.397518943312669
$arr_B_E[23,2] = 15.02774251200988
$arr_B_E[23,3] = 16.43910790410945
$ws.Range("B2:E25").Value = $arr_B_E

$arr_G_K = New-Object 'object[,]' 24,5
$arr_G_K[0,0] = 32.91174236810749
$arr_G_K[0,1] = 15.388259495514
$arr_G_K[0,2] = 21.87708099513625
$arr_G_K[0,3] = 9.317497202152996
$arr_G_K[0,4] = 10.23282874145926
$arr_G_K[1,0] = 32.96866988064129
$arr_G_K[1,1] = 15.43241214227752
$arr_G_K[1,2] = 21.95943644551062
$arr_G_K[1,3] = 9.323620176607356
$arr_G_K[1,4] = 10.04500976814862
$arr_G_K[2,0] = 33.01278252662116
$arr_G_K[2,1] = 15.4617916996163
$arr_G_K[2,2] = 22.01380304982681
$arr_G_K[2,3] = 9.328757147649059
$arr_G_K[2,4] = 9.929409786853464
$arr_G_K[3,0] = 33.03305548491772
$arr_G_K[3,1] = 15.47433486047323
$arr_G_K[3,2] = 22.03691304031437
$arr_G_K[3,3] = 9.331197257912512
$arr_G_K[3,4] = 9.882293010168725
$arr_G_K[4,0] = 33.03656030398368
$arr_G_K[4,1] = 15.4764521156716
$arr_G_K[4,2] = 22.04080811477826
$arr_G_K[4,3] = 9.331623388962882
$arr_G_K[4,4] = 9.874470579985433
$arr_G_K[5,0] = 33.01304664474727
$arr_G_K[5,1] = 15.46195855008302
$arr_G_K[5,2] = 22.01411085235263
$arr_G_K[5,3] = 9.32878865139398
$arr_G_K[5,4] = 9.928774307722152
$arr_G_K[6,0] = 32.92946647517102
$arr_G_K[6,1] = 15.4030122817795
$arr_G_K[6,2] = 21.90468822475396
$arr_G_K[6,3] = 9.319322769702348
$arr_G_K[6,4] = 10.16816284403659
$arr_G_K[7,0] = 32.83848351122342
$arr_G_K[7,1] = 15.30542610811688
$arr_G_K[7,2] = 21.72027965184807
$arr_G_K[7,3] = 9.311670229358558
$arr_G_K[7,4] = 10.63277759677632
$arr_G_K[8,0] = 32.81636621987104
$arr_G_K[8,1] = 15.24470519521818
$arr_G_K[8,2] = 21.60321016809769
$arr_G_K[8,3] = 9.31266928418369
$arr_G_K[8,4] = 10.96781188586938
$arr_G_K[9,0] = 32.81605361501067
$arr_G_K[9,1] = 15.21946447663606
$arr_G_K[9,2] = 21.5539564976975
$arr_G_K[9,3] = 9.314553422553832
$arr_G_K[9,4] = 11.11818544211532
$arr_G_K[10,0] = 32.81733816542675
$arr_G_K[10,1] = 15.21024892414096
$arr_G_K[10,2] = 21.53588144795956
$arr_G_K[10,3] = 9.315471609017031
$arr_G_K[10,4] = 11.17478297907072
$arr_G_K[11,0] = 32.8169991143445
$arr_G_K[11,1] = 15.21221842132752
$arr_G_K[11,2] = 21.53974859189479
$arr_G_K[11,3] = 9.315264771920523
$arr_G_K[11,4] = 11.16260987340334
$arr_G_K[12,0] = 32.81613117701595
$arr_G_K[12,1] = 15.21869944163893
$arr_G_K[12,2] = 21.55245789907177
$arr_G_K[12,3] = 9.314624865120196
$arr_G_K[12,4] = 11.12284897646481
$arr_G_K[13,0] = 32.81578225474438
$arr_G_K[13,1] = 15.22271386791612
$arr_G_K[13,2] = 21.56031778463665
$arr_G_K[13,3] = 9.314259534458841
$arr_G_K[13,4] = 11.09844773960898
$arr_G_K[14,0] = 32.8165829367459
$arr_G_K[14,1] = 15.24640266661995
$arr_G_K[14,2] = 21.6065096053142
$arr_G_K[14,3] = 9.312574846846996
$arr_G_K[14,4] = 10.95793906987756
$arr_G_K[15,0] = 32.81957213907314
$arr_G_K[15,1] = 15.26154501687323
$arr_G_K[15,2] = 21.63587231082714
$arr_G_K[15,3] = 9.311906922062184
$arr_G_K[15,4] = 10.87118166230462
$arr_G_K[16,0] = 32.82220904894746
$arr_G_K[16,1] = 15.27047865156752
$arr_G_K[16,2] = 21.65313758492702
$arr_G_K[16,3] = 9.31165737494848
$arr_G_K[16,4] = 10.82109252342106
$arr_G_K[17,0] = 32.82325940772255
$arr_G_K[17,1] = 15.2735419225426
$arr_G_K[17,2] = 21.65904797159326
$arr_G_K[17,3] = 9.311596029840508
$arr_G_K[17,4] = 10.80410247614944
$arr_G_K[18,0] = 32.81915895364497
$arr_G_K[18,1] = 15.25990988785463
$arr_G_K[18,2] = 21.6327076160078
$arr_G_K[18,3] = 9.311964096616194
$arr_G_K[18,4] = 10.88043704610152
$arr_G_K[19,0] = 32.81634803346488
$arr_G_K[19,1] = 15.21678651114702
$arr_G_K[19,2] = 21.54870922195528
$arr_G_K[19,3] = 9.314807273147386
$arr_G_K[19,4] = 11.13453750498122
$arr_G_K[20,0] = 32.82268792857347
$arr_G_K[20,1] = 15.1905996325901
$arr_G_K[20,2] = 21.49717059661572
$arr_G_K[20,3] = 9.317858147985419
$arr_G_K[20,4] = 11.29856952732426
$arr_G_K[21,0] = 32.81855596155735
$arr_G_K[21,1] = 15.20439333001415
$arr_G_K[21,2] = 21.52437010116443
$arr_G_K[21,3] = 9.31612102317475
$arr_G_K[21,4] = 11.21122571992786
$arr_G_K[22,0] = 32.81934289417692
$arr_G_K[22,1] = 15.26064841917718
$arr_G_K[22,2] = 21.63413717762664
$arr_G_K[22,3] = 9.311937829133104
$arr_G_K[22,4] = 10.87625334069832
$arr_G_K[23,0] = 32.8552578474067
$arr_G_K[23,1] = 15.32989814673271
$arr_G_K[23,2] = 21.76693489017619
$arr_G_K[23,3] = 9.312575144088768
$arr_G_K[23,4] = 10.50796053440034
$ws.Range("G2:K25").Value = $arr_G_K

$arr_N_O = New-Object 'object[,]' 24,2
$arr_N_O[0,0] = 18.69568457314384
$arr_N_O[0,1] = 23.96746889359511
$arr_N_O[1,0] = 18.75115681069127
$arr_N_O[1,1] = 24.03537511205979
$arr_N_O[2,0] = 18.78683530315897
$arr_N_O[2,1] = 24.08171960257835
$arr_N_O[3,0] = 18.80178268929337
$arr_N_O[3,1] = 24.10177286419333
$arr_N_O[4,0] = 18.80428937479256
$arr_N_O[4,1] = 24.105173155234
$arr_N_O[5,0] = 18.78703523485687
$arr_N_O[5,1] = 24.08198532305257
$arr_N_O[6,0] = 18.71447622799216
$arr_N_O[6,1] = 23.98991688348269
$arr_N_O[7,0] = 18.58497452514754
$arr_N_O[7,1] = 23.84633912240713
$arr_N_O[8,0] = 18.4975473203264
$arr_N_O[8,1] = 23.7634861021739
$arr_N_O[9,0] = 18.45943381449056
$arr_N_O[9,1] = 23.73072737725036
$arr_N_O[10,0] = 18.44523836853336
$arr_N_O[10,1] = 23.71903298470383
$arr_N_O[11,0] = 18.44828507382694
$arr_N_O[11,1] = 23.72151995403025
$arr_N_O[12,0] = 18.45826119647499
$arr_N_O[12,1] = 23.72975102256089
$arr_N_O[13,0] = 18.46440273380902
$arr_N_O[13,1] = 23.73488537567559
$arr_N_O[14,0] = 18.50007139250152
$arr_N_O[14,1] = 23.76572632626455
$arr_N_O[15,0] = 18.52237675521925
$arr_N_O[15,1] = 23.78591042750286
$arr_N_O[16,0] = 18.53536228975367
$arr_N_O[16,1] = 23.79798388752479
$arr_N_O[17,0] = 18.5397858075248
$arr_N_O[17,1] = 23.80215141585596
$arr_N_O[18,0] = 18.5199861651976
$arr_N_O[18,1] = 23.78371375183509
$arr_N_O[19,0] = 18.45532453512981
$arr_N_O[19,1] = 23.72731406093087
$arr_N_O[20,0] = 18.41444717857633
$arr_N_O[20,1] = 23.69459626096873
$arr_N_O[21,0] = 18.43613801293704
$arr_N_O[21,1] = 23.71167886023843
$arr_N_O[22,0] = 18.52106644659625
$arr_N_O[22,1] = 23.78470540700488
$arr_N_O[23,0] = 18.61864732662093
$arr_N_O[23,1] = 23.88121257723607
$ws.Range("N2:O25").Value = $arr_N_O

